$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("FACTORS_MUNI")
$ws1.Name = "FACTORS"
